# Update LR-pair data rows with refreshed TPM-based values.
# Adds a new "Resolving-Mac" sending-cluster block (rows 17-21) and refreshes
# every existing numeric column (E:T) plus the ligand/receptor index pairs (B:C)
# for the previously-present rows (2-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @(2, "ECs", "Npnt", "Itga8", "ECs", 3, 1, 2.128539, 6.385617, 0.4682231486221785, 0.4682231486221785, 3, 1, 0.3708513333333334, 1.112554, 0.07446453063452614, 0.07446453063452615, 0.789371526202, 7.104343735818, 0.0348660169943705, 0.0348660169943705),
    @(3, "ECs", "Npnt", "Itga8", "FAPs", 3, 1, 2.128539, 6.385617, 0.4682231486221785, 0.4682231486221785, 3, 1, 2.230730666666667, 6.692192, 0.4479161786269527, 0.4479161786269528, 4.748197222496, 42.733775002464, 0.2097247234755259, 0.2097247234755259),
    @(4, "ECs", "Npnt", "Itga8", "Inflammatory-Mac", 3, 1, 2.128539, 6.385617, 0.4682231486221785, 0.4682231486221785, 3, 1, 1.029162, 3.087486, 0.2066490218278579, 0.2066490218278579, 2.190611454318, 19.715503088862, 0.0967578556599329, 0.0967578556599329),
    @(5, "ECs", "Npnt", "Itga8", "MuSCs", 3, 1, 2.128539, 6.385617, 0.4682231486221785, 0.4682231486221785, 3, 1, 0.8350726666666667, 2.505218, 0.1676771487111334, 0.1676771487111334, 1.777484738834, 15.997362649506, 0.07851032252151613, 0.07851032252151614),
    @(6, "ECs", "Npnt", "Itga8", "Resolving-Mac", 3, 1, 2.128539, 6.385617, 0.4682231486221785, 0.4682231486221785, 3, 1, 0.5144246666666666, 1.543274, 0.1032931201995298, 0.1032931201995298, 1.094972965562, 9.854756690058, 0.04836422997083299, 0.048364229970833),
    @(7, "FAPs", "Npnt", "Itga8", "ECs", 3, 1, 0.6061233333333333, 1.81837, 0.133331348679401, 0.133331348679401, 3, 1, 0.3708513333333334, 1.112554, 0.07446453063452614, 0.07446453063452615, 0.2247816463311111, 2.02303481698, 0.009928456298279943, 0.009928456298279945),
    @(8, "FAPs", "Npnt", "Itga8", "FAPs", 3, 1, 0.6061233333333333, 1.81837, 0.133331348679401, 0.133331348679401, 3, 1, 2.230730666666667, 6.692192, 0.4479161786269527, 0.4479161786269528, 1.352097907448889, 12.16888116704, 0.05972126819165509, 0.0597212681916551),
    @(9, "FAPs", "Npnt", "Itga8", "Inflammatory-Mac", 3, 1, 0.6061233333333333, 1.81837, 0.133331348679401, 0.133331348679401, 3, 1, 1.029162, 3.087486, 0.2066490218278579, 0.2066490218278579, 0.6237991019800001, 5.61419191782, 0.02755279278358727, 0.02755279278358727),
    @(10, "FAPs", "Npnt", "Itga8", "MuSCs", 3, 1, 0.6061233333333333, 1.81837, 0.133331348679401, 0.133331348679401, 3, 1, 0.8350726666666667, 2.505218, 0.1676771487111334, 0.1676771487111334, 0.5061570282955555, 4.55541325466, 0.0223566203803719, 0.0223566203803719),
    @(11, "FAPs", "Npnt", "Itga8", "Resolving-Mac", 3, 1, 0.6061233333333333, 1.81837, 0.133331348679401, 0.133331348679401, 3, 1, 0.5144246666666666, 1.543274, 0.1032931201995298, 0.1032931201995298, 0.3118047937088889, 2.80624314338, 0.01377221102550679, 0.01377221102550679),
    @(12, "MuSCs", "Npnt", "Itga8", "ECs", 3, 1, 1.392473666666667, 4.177421, 0.306307943890216, 0.306307943890216, 3, 1, 0.3708513333333334, 1.112554, 0.07446453063452614, 0.07446453063452615, 0.5164007159148889, 4.647606443234, 0.0228090772714117, 0.02280907727141171),
    @(13, "MuSCs", "Npnt", "Itga8", "FAPs", 3, 1, 1.392473666666667, 4.177421, 0.306307943890216, 0.306307943890216, 3, 1, 2.230730666666667, 6.692192, 0.4479161786269527, 0.4479161786269528, 3.106233710759111, 27.956103396832, 0.1372002837103846, 0.1372002837103846),
    @(14, "MuSCs", "Npnt", "Itga8", "Inflammatory-Mac", 3, 1, 1.392473666666667, 4.177421, 0.306307943890216, 0.306307943890216, 3, 1, 1.029162, 3.087486, 0.2066490218278579, 0.2066490218278579, 1.433080983734, 12.897728853606, 0.06329823698301551, 0.06329823698301551),
    @(15, "MuSCs", "Npnt", "Itga8", "MuSCs", 3, 1, 1.392473666666667, 4.177421, 0.306307943890216, 0.306307943890216, 3, 1, 0.8350726666666667, 2.505218, 0.1676771487111334, 0.1676771487111334, 1.162816698086445, 10.465350282778, 0.05136084265908125, 0.05136084265908126),
    @(16, "MuSCs", "Npnt", "Itga8", "Resolving-Mac", 3, 1, 1.392473666666667, 4.177421, 0.306307943890216, 0.306307943890216, 3, 1, 0.5144246666666666, 1.543274, 0.1032931201995298, 0.1032931201995298, 0.716322801817111, 6.446905216354, 0.03163950326632291, 0.03163950326632291),
    @(17, "Resolving-Mac", "Npnt", "Itga8", "ECs", 2, 0.6666666666666666, 0.4188566666666667, 1.25657, 0.09213755880820455, 0.09213755880820455, 3, 1, 0.3708513333333334, 1.112554, 0.07446453063452614, 0.07446453063452615, 0.1553335533088889, 1.39800197978, 0.006860980070464002, 0.006860980070464003),
    @(18, "Resolving-Mac", "Npnt", "Itga8", "FAPs", 2, 0.6666666666666666, 0.4188566666666667, 1.25657, 0.09213755880820455, 0.09213755880820455, 3, 1, 2.230730666666667, 6.692192, 0.4479161786269527, 0.4479161786269528, 0.934356411271111, 8.40920770144, 0.04126990324938711, 0.04126990324938712),
    @(19, "Resolving-Mac", "Npnt", "Itga8", "Inflammatory-Mac", 2, 0.6666666666666666, 0.4188566666666667, 1.25657, 0.09213755880820455, 0.09213755880820455, 3, 1, 1.029162, 3.087486, 0.2066490218278579, 0.2066490218278579, 0.43107136478, 3.87964228302, 0.0190401364013222, 0.0190401364013222),
    @(20, "Resolving-Mac", "Npnt", "Itga8", "MuSCs", 2, 0.6666666666666666, 0.4188566666666667, 1.25657, 0.09213755880820455, 0.09213755880820455, 3, 1, 0.8350726666666667, 2.505218, 0.1676771487111334, 0.1676771487111334, 0.3497757535844445, 3.14798178226, 0.01544936315016411, 0.01544936315016412),
    @(21, "Resolving-Mac", "Npnt", "Itga8", "Resolving-Mac", 2, 0.6666666666666666, 0.4188566666666667, 1.25657, 0.09213755880820455, 0.09213755880820455, 3, 1, 0.5144246666666666, 1.543274, 0.1032931201995298, 0.1032931201995298, 0.2154702011311111, 1.93923181018, 0.00951717593686712, 0.009517175936867121)
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
    $ws.Cells.Item($r, 18).Value = $row[18]
    $ws.Cells.Item($r, 19).Value = $row[19]
    $ws.Cells.Item($r, 20).Value = $row[20]
}

